# Reconciliation of capital and lowercase letters in the "Category" column
# (column H) of the "Timeseries Attributes" tab, so the values match the
# capitalization used for the visualization parameters elsewhere
# (load -> Load, solar -> Solar, wind -> Wind).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Timeseries Attributes")

# Make this the active tab, as it was when the edit was made.
$ws.Activate()

$ws.Range("H2").Value = "Load"
$ws.Range("H3").Value = "Load"
$ws.Range("H4").Value = "Solar"
$ws.Range("H5").Value = "Wind"

# Leave the selection on the last cell that was typed into.
$ws.Range("H5").Select()
